# Scheduled-runner style update: refresh market-price-derived columns
# (H currentAveragePrice, I currentAveragePriceNQ, J currentAveragePriceHQ,
#  K LevePriceNQ, L LevePriceHQ, M LeveProfitNQ, N LeveProfitHQ) across the
# leve-profit tables on each crafting-job sheet. Cells are plain cached
# numeric values (no formulas in this workbook), so each changed cell is
# written directly; a few rows gain/lose trailing profit cells entirely.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 76 - Warding Off Temptation
$ws.Cells.Item(76, 8).Value = 7268.0967
$ws.Cells.Item(76, 9).Value = 12382.091
$ws.Cells.Item(76, 10).Value = 4455.4
$ws.Cells.Item(76, 11).Value = 12382.091
$ws.Cells.Item(76, 12).Value = 4455.4
$ws.Cells.Item(76, 13).Value = -12067.091
$ws.Cells.Item(76, 14).Value = -5085.4

# ALC row 79 - The Garden of Arcane Delights (L)
$ws.Cells.Item(79, 8).Value = 7268.0967
$ws.Cells.Item(79, 9).Value = 12382.091
$ws.Cells.Item(79, 10).Value = 4455.4
$ws.Cells.Item(79, 11).Value = 12382.091
$ws.Cells.Item(79, 12).Value = 4455.4
$ws.Cells.Item(79, 13).Value = -11290.091
$ws.Cells.Item(79, 14).Value = -6639.4

# ALC row 92 - Whinier than the Sword
$ws.Cells.Item(92, 8).Value = 73099860
$ws.Cells.Item(92, 9).Value = 3472711.5
$ws.Cells.Item(92, 10).Value = 444444600
$ws.Cells.Item(92, 11).Value = 3472711.5
$ws.Cells.Item(92, 12).Value = 444444600
$ws.Cells.Item(92, 13).Value = -3471463.5
$ws.Cells.Item(92, 14).Value = -444447096

# ALC row 137 - Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 3986.5
$ws.Cells.Item(137, 9).Value = 2549.4
$ws.Cells.Item(137, 10).Value = 6381.6665
$ws.Cells.Item(137, 11).Value = 7648.200000000001
$ws.Cells.Item(137, 12).Value = 19144.9995
$ws.Cells.Item(137, 13).Value = -5098.200000000001
$ws.Cells.Item(137, 14).Value = -24244.9995

# ALC row 138 - All-night Crafting
$ws.Cells.Item(138, 8).Value = 1936.9697
$ws.Cells.Item(138, 9).Value = 865.913
$ws.Cells.Item(138, 10).Value = 2866.566
$ws.Cells.Item(138, 11).Value = 2597.739
$ws.Cells.Item(138, 12).Value = 8599.698
$ws.Cells.Item(138, 13).Value = 2542.261
$ws.Cells.Item(138, 14).Value = -18879.698

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45 - Hollow Hallmarks
$ws.Cells.Item(45, 8).Value = 8225.875
$ws.Cells.Item(45, 9).Value = 9937.23
$ws.Cells.Item(45, 10).Value = 810
$ws.Cells.Item(45, 11).Value = 9937.23
$ws.Cells.Item(45, 12).Value = 810
$ws.Cells.Item(45, 13).Value = -9560.23
$ws.Cells.Item(45, 14).Value = -1564

# ARM row 61 - Dealing with the Tough Stuff
$ws.Cells.Item(61, 8).Value = 249125.73
$ws.Cells.Item(61, 9).Value = 6119.9614
$ws.Cells.Item(61, 10).Value = 670335.75
$ws.Cells.Item(61, 11).Value = 6119.9614
$ws.Cells.Item(61, 12).Value = 670335.75
$ws.Cells.Item(61, 13).Value = -5907.9614
$ws.Cells.Item(61, 14).Value = -670759.75

# ARM row 63 - Rivets Run through It
$ws.Cells.Item(63, 8).Value = 111124104
$ws.Cells.Item(63, 9).Value = 125014240
$ws.Cells.Item(63, 11).Value = 125014240
$ws.Cells.Item(63, 13).Value = -125013554

# ARM row 66 - A Riveting Revival (L)
$ws.Cells.Item(66, 8).Value = 111124104
$ws.Cells.Item(66, 9).Value = 125014240
$ws.Cells.Item(66, 11).Value = 625071200
$ws.Cells.Item(66, 13).Value = -625067768

# ARM row 74 - As the Bolt Flies
$ws.Cells.Item(74, 8).Value = 2390.3684
$ws.Cells.Item(74, 9).Value = 2083.5833
$ws.Cells.Item(74, 10).Value = 2916.2856
$ws.Cells.Item(74, 11).Value = 2083.5833
$ws.Cells.Item(74, 12).Value = 2916.2856
$ws.Cells.Item(74, 13).Value = -1209.5833
$ws.Cells.Item(74, 14).Value = -4664.2856

# ARM row 77 - Heavy Metal Banned (L)
$ws.Cells.Item(77, 8).Value = 2390.3684
$ws.Cells.Item(77, 9).Value = 2083.5833
$ws.Cells.Item(77, 10).Value = 2916.2856
$ws.Cells.Item(77, 11).Value = 10417.9165
$ws.Cells.Item(77, 12).Value = 14581.428
$ws.Cells.Item(77, 13).Value = -6049.916499999999
$ws.Cells.Item(77, 14).Value = -23317.428

# ARM row 110 - Scheduled Maintenance
$ws.Cells.Item(110, 8).Value = 1231.2413
$ws.Cells.Item(110, 9).Value = 1201.8695
$ws.Cells.Item(110, 10).Value = 1343.8334
$ws.Cells.Item(110, 11).Value = 1201.8695
$ws.Cells.Item(110, 12).Value = 1343.8334
$ws.Cells.Item(110, 13).Value = 843.1305
$ws.Cells.Item(110, 14).Value = -5433.8334

# ARM row 136 - Metal with Mettle
$ws.Cells.Item(136, 8).Value = 249125.73
$ws.Cells.Item(136, 9).Value = 6119.9614
$ws.Cells.Item(136, 10).Value = 670335.75
$ws.Cells.Item(136, 11).Value = 18359.8842
$ws.Cells.Item(136, 12).Value = 2011007.25
$ws.Cells.Item(136, 13).Value = -15809.8842
$ws.Cells.Item(136, 14).Value = -2016107.25

$ws = $wb.Worksheets.Item("BSM")
# BSM row 35 - Lancers' Creed
$ws.Cells.Item(35, 8).Value = 21000
$ws.Cells.Item(35, 9).Value = 20000
$ws.Cells.Item(35, 10).Value = 22000
$ws.Cells.Item(35, 11).Value = 20000
$ws.Cells.Item(35, 12).Value = 22000
$ws.Cells.Item(35, 13).Value = -19690
$ws.Cells.Item(35, 14).Value = -22620

# BSM row 82 - Spirituality Inspector
$ws.Cells.Item(82, 8).Value = 11882.8
$ws.Cells.Item(82, 10).Value = 25466.666
$ws.Cells.Item(82, 12).Value = 25466.666
$ws.Cells.Item(82, 14).Value = -26232.666

# BSM row 85 - The Clamor for Hammers (L)
$ws.Cells.Item(85, 8).Value = 11882.8
$ws.Cells.Item(85, 10).Value = 25466.666
$ws.Cells.Item(85, 12).Value = 25466.666
$ws.Cells.Item(85, 14).Value = -28118.666

# BSM row 86 - Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 2004.5
$ws.Cells.Item(86, 9).Value = 1938.3334
$ws.Cells.Item(86, 10).Value = 2600
$ws.Cells.Item(86, 11).Value = 1938.3334
$ws.Cells.Item(86, 12).Value = 2600
$ws.Cells.Item(86, 13).Value = -815.3334
$ws.Cells.Item(86, 14).Value = -4846

# BSM row 89 - Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 2004.5
$ws.Cells.Item(89, 9).Value = 1938.3334
$ws.Cells.Item(89, 10).Value = 2600
$ws.Cells.Item(89, 11).Value = 9691.666999999999
$ws.Cells.Item(89, 12).Value = 13000
$ws.Cells.Item(89, 13).Value = -4075.666999999999
$ws.Cells.Item(89, 14).Value = -24232

# BSM row 99 - Meddle in Metal
$ws.Cells.Item(99, 8).Value = 83335330
$ws.Cells.Item(99, 9).Value = 166667920
$ws.Cells.Item(99, 10).Value = 2739.1667
$ws.Cells.Item(99, 11).Value = 166667920
$ws.Cells.Item(99, 12).Value = 2739.1667
$ws.Cells.Item(99, 13).Value = -166666422
$ws.Cells.Item(99, 14).Value = -5735.1667

$ws = $wb.Worksheets.Item("CRP")
# CRP row 4 - A Clogful of Camaraderie
$ws.Cells.Item(4, 8).Value = 100002
$ws.Cells.Item(4, 10).Value = 100002
$ws.Cells.Item(4, 12).Value = 100002
$ws.Cells.Item(4, 14).Value = -100226

# CRP row 58 - You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 251620.9
$ws.Cells.Item(58, 9).Value = 1406.625
$ws.Cells.Item(58, 10).Value = 626942.3
$ws.Cells.Item(58, 11).Value = 1406.625
$ws.Cells.Item(58, 12).Value = 626942.3
$ws.Cells.Item(58, 13).Value = -1203.625
$ws.Cells.Item(58, 14).Value = -627348.3

# CRP row 107 - Built to Last
$ws.Cells.Item(107, 8).Value = 836.62964
$ws.Cells.Item(107, 9).Value = 656.41174
$ws.Cells.Item(107, 10).Value = 1143
$ws.Cells.Item(107, 11).Value = 656.41174
$ws.Cells.Item(107, 12).Value = 1143
$ws.Cells.Item(107, 13).Value = 1263.58826
$ws.Cells.Item(107, 14).Value = -4983

# CRP row 108 - Just Starting Out
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 9).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 11).Value = 0
$ws.Cells.Item(108, 12).ClearContents()  # L108: was 40000, now blank
$ws.Cells.Item(108, 13).ClearContents()  # M108: was -970.5, now blank
$ws.Cells.Item(108, 14).Value = 0

# CRP row 109 - Playing the Market
$ws.Cells.Item(109, 8).Value = 42942.5
$ws.Cells.Item(109, 10).Value = 42942.5
$ws.Cells.Item(109, 12).Value = 42942.5
$ws.Cells.Item(109, 14).Value = -45022.5

# CRP row 110 - A Stronger Offense
$ws.Cells.Item(110, 8).Value = 43000
$ws.Cells.Item(110, 10).Value = 43000
$ws.Cells.Item(110, 12).Value = 43000
$ws.Cells.Item(110, 14).Value = -51180

# CRP row 122 - Timber of Tenkonto
$ws.Cells.Item(122, 8).Value = 1502.4166
$ws.Cells.Item(122, 9).Value = 1502.4166
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 4507.2498
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).ClearContents()  # M122: was -2247.625, now blank
$ws.Cells.Item(122, 14).Value = -2057.2498

# CRP row 136 - Turali Quality
$ws.Cells.Item(136, 8).Value = 251620.9
$ws.Cells.Item(136, 9).Value = 1406.625
$ws.Cells.Item(136, 10).Value = 626942.3
$ws.Cells.Item(136, 11).Value = 4219.875
$ws.Cells.Item(136, 12).Value = 1880826.9
$ws.Cells.Item(136, 13).Value = -1669.875
$ws.Cells.Item(136, 14).Value = -1885926.9

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5 - What a Sap
$ws.Cells.Item(5, 8).Value = 9844.846
$ws.Cells.Item(5, 10).Value = 3485.7144
$ws.Cells.Item(5, 12).Value = 10457.1432
$ws.Cells.Item(5, 14).Value = -10681.1432

# CUL row 131 - The Mountain Steeped
$ws.Cells.Item(131, 8).Value = 1755323.6
$ws.Cells.Item(131, 9).Value = 3846553.8
$ws.Cells.Item(131, 10).Value = 1388.7097
$ws.Cells.Item(131, 11).Value = 11539661.4
$ws.Cells.Item(131, 12).Value = 4166.1291
$ws.Cells.Item(131, 13).Value = -11534621.4
$ws.Cells.Item(131, 14).Value = -14246.1291

# CUL row 135 - Not-so-secret Ingredient
$ws.Cells.Item(135, 8).Value = 9844.846
$ws.Cells.Item(135, 10).Value = 3485.7144
$ws.Cells.Item(135, 12).Value = 31371.4296
$ws.Cells.Item(135, 14).Value = -36441.4296

$ws = $wb.Worksheets.Item("GSM")
# GSM row 15 - The Tusk at Hand
$ws.Cells.Item(15, 8).Value = 30000
$ws.Cells.Item(15, 10).Value = 30000
$ws.Cells.Item(15, 12).Value = 30000
$ws.Cells.Item(15, 14).Value = -30576

# GSM row 80 - Needs More Prayerbell
$ws.Cells.Item(80, 8).Value = 9457.5
$ws.Cells.Item(80, 9).Value = 13322.777
$ws.Cells.Item(80, 10).Value = 2500
$ws.Cells.Item(80, 11).Value = 13322.777
$ws.Cells.Item(80, 12).Value = 2500
$ws.Cells.Item(80, 13).Value = -12324.777
$ws.Cells.Item(80, 14).Value = -4496

# GSM row 81 - The Grander Temple
$ws.Cells.Item(81, 8).Value = 30000
$ws.Cells.Item(81, 10).Value = 30000
$ws.Cells.Item(81, 12).Value = 30000
$ws.Cells.Item(81, 14).Value = -31996

# GSM row 83 - With a Noise That Reaches Heaven (L)
$ws.Cells.Item(83, 8).Value = 9457.5
$ws.Cells.Item(83, 9).Value = 13322.777
$ws.Cells.Item(83, 10).Value = 2500
$ws.Cells.Item(83, 11).Value = 66613.88499999999
$ws.Cells.Item(83, 12).Value = 12500
$ws.Cells.Item(83, 13).Value = -61621.88499999999
$ws.Cells.Item(83, 14).Value = -22484

# GSM row 84 - Man with a Dragon Earring (L)
$ws.Cells.Item(84, 8).Value = 30000
$ws.Cells.Item(84, 10).Value = 30000
$ws.Cells.Item(84, 12).Value = 90000
$ws.Cells.Item(84, 14).Value = -99984

# GSM row 122 - Awarding Academic Excellence
$ws.Cells.Item(122, 8).Value = 8747682
$ws.Cells.Item(122, 9).Value = 3242801
$ws.Cells.Item(122, 10).Value = 15628784
$ws.Cells.Item(122, 11).Value = 9728403
$ws.Cells.Item(122, 12).Value = 46886352
$ws.Cells.Item(122, 13).Value = -9725953
$ws.Cells.Item(122, 14).Value = -46891252

$ws = $wb.Worksheets.Item("LTW")
# LTW row 132 - Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 8338839
$ws.Cells.Item(132, 9).Value = 10422840
$ws.Cells.Item(132, 10).Value = 2837.375
$ws.Cells.Item(132, 11).Value = 31268520
$ws.Cells.Item(132, 12).Value = 8512.125
$ws.Cells.Item(132, 13).Value = -31265990
$ws.Cells.Item(132, 14).Value = -13572.125

$ws = $wb.Worksheets.Item("WVR")
# WVR row 15 - Workplace Safety
$ws.Cells.Item(15, 8).Value = 7200
$ws.Cells.Item(15, 10).Value = 7200
$ws.Cells.Item(15, 12).Value = 7200
$ws.Cells.Item(15, 14).Value = -7776

# WVR row 19 - Dirt Cheap
$ws.Cells.Item(19, 8).Value = 1050
$ws.Cells.Item(19, 9).Value = 100
$ws.Cells.Item(19, 10).Value = 2000
$ws.Cells.Item(19, 11).Value = 100
$ws.Cells.Item(19, 12).Value = 2000
$ws.Cells.Item(19, 13).Value = 74
$ws.Cells.Item(19, 14).Value = -2348

# WVR row 81 - Where the Dragonflies, the Net Catches
$ws.Cells.Item(81, 8).Value = 1589
$ws.Cells.Item(81, 9).Value = 1157.2858
$ws.Cells.Item(81, 10).Value = 3100
$ws.Cells.Item(81, 11).Value = 2314.5716
$ws.Cells.Item(81, 12).Value = 6200
$ws.Cells.Item(81, 13).Value = -1253.5716
$ws.Cells.Item(81, 14).Value = -8322

# WVR row 84 - To Kill a Dragon on Nameday (L)
$ws.Cells.Item(84, 8).Value = 1589
$ws.Cells.Item(84, 9).Value = 1157.2858
$ws.Cells.Item(84, 10).Value = 3100
$ws.Cells.Item(84, 11).Value = 11572.858
$ws.Cells.Item(84, 12).Value = 31000
$ws.Cells.Item(84, 13).Value = -6268.858
$ws.Cells.Item(84, 14).Value = -41608

# WVR row 94 - Proper Props
$ws.Cells.Item(94, 8).Value = 27311.8
$ws.Cells.Item(94, 10).Value = 27311.8
$ws.Cells.Item(94, 12).Value = 27311.8
$ws.Cells.Item(94, 14).Value = -29113.8

# WVR row 132 - Comfy Cabins
$ws.Cells.Item(132, 8).Value = 1127.5062
$ws.Cells.Item(132, 9).Value = 816.9231
$ws.Cells.Item(132, 10).Value = 2389.25
$ws.Cells.Item(132, 11).Value = 2450.7693
$ws.Cells.Item(132, 12).Value = 7167.75
$ws.Cells.Item(132, 13).Value = 79.23070000000007
$ws.Cells.Item(132, 14).Value = -12227.75
